$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("p1.xlsx")

# Row 7
$ws.Range("A7").Value = 0.0
$ws.Range("D7").Value = 469.0
$ws.Range("F7").Value = 4.0
$ws.Range("G7").Value = 14.0
$ws.Range("H7").Value = 6.0
$ws.Range("I7").Value = 27.0
$ws.Range("J7").Value = 51.0
$ws.Range("K7").Value = 32.0
$ws.Range("L7").Value = 26.0
$ws.Range("M7").Value = 22.0
$ws.Range("N7").Value = 28.0
$ws.Range("O7").Value = 3.0
$ws.Range("P7").Value = 36.0
$ws.Range("Q7").Value = 35.0
$ws.Range("R7").Value = 20.0
$ws.Range("S7").Value = 2.0
$ws.Range("T7").Value = 50.0
$ws.Range("U7").Value = 9.0
$ws.Range("V7").Value = 10.0
$ws.Range("W7").Value = 5.0
$ws.Range("X7").Value = 12.0
$ws.Range("Y7").Value = 17.0
$ws.Range("Z7").Value = 37.0
$ws.Range("AA7").Value = 15.0
$ws.Range("AB7").Value = 42.0
$ws.Range("AC7").Value = -1.0
$ws.Range("AD7").ClearContents()
$ws.Range("AE7").ClearContents()

# Row 8
$ws.Range("A8").Value = 0.0
$ws.Range("D8").Value = 468.0
$ws.Range("F8").Value = 41.0
$ws.Range("G8").Value = 13.0
$ws.Range("H8").Value = 47.0
$ws.Range("I8").Value = 18.0
$ws.Range("J8").Value = 25.0
$ws.Range("K8").Value = 24.0
$ws.Range("L8").Value = 43.0
$ws.Range("M8").Value = 23.0
$ws.Range("N8").Value = 7.0
$ws.Range("O8").Value = 48.0
$ws.Range("P8").Value = 8.0
$ws.Range("Q8").Value = 31.0
$ws.Range("R8").Value = 1.0
$ws.Range("S8").Value = 29.0
$ws.Range("T8").Value = 21.0
$ws.Range("U8").Value = 34.0
$ws.Range("V8").Value = 30.0
$ws.Range("W8").Value = 16.0
$ws.Range("Y8").Value = 46.0
$ws.Range("Z8").Value = 38.0
$ws.Range("AA8").Value = 49.0
$ws.Range("AB8").Value = 39.0
$ws.Range("AC8").Value = 33.0
$ws.Range("AD8").Value = 45.0
$ws.Range("AE8").Value = 44.0
$ws.Range("AF8").Value = 19.0
$ws.Range("AG8").Value = 40.0
$ws.Range("AH8").Value = -1.0

# Row 9
$ws.Range("A9").Value = 0.0
$ws.Range("D9").Value = 469.0
$ws.Range("F9").Value = 4.0
$ws.Range("G9").Value = 14.0
$ws.Range("H9").Value = 6.0
$ws.Range("I9").Value = 27.0
$ws.Range("J9").Value = 51.0
$ws.Range("K9").Value = 32.0
$ws.Range("L9").Value = 26.0
$ws.Range("M9").Value = 22.0
$ws.Range("N9").Value = 28.0
$ws.Range("O9").Value = 3.0
$ws.Range("P9").Value = 36.0
$ws.Range("Q9").Value = 35.0
$ws.Range("R9").Value = 20.0
$ws.Range("S9").Value = 2.0
$ws.Range("T9").Value = 50.0
$ws.Range("U9").Value = 9.0
$ws.Range("V9").Value = 10.0
$ws.Range("W9").Value = 5.0
$ws.Range("X9").Value = 12.0
$ws.Range("Y9").Value = 17.0
$ws.Range("Z9").Value = 37.0
$ws.Range("AA9").Value = 15.0
$ws.Range("AB9").Value = 42.0
$ws.Range("AC9").Value = -1.0
$ws.Range("AD9").ClearContents()
$ws.Range("AE9").ClearContents()

# Row 10
$ws.Range("A10").Value = 0.0
$ws.Range("D10").Value = 468.0
$ws.Range("F10").Value = 41.0
$ws.Range("G10").Value = 13.0
$ws.Range("H10").Value = 47.0
$ws.Range("I10").Value = 18.0
$ws.Range("J10").Value = 25.0
$ws.Range("K10").Value = 24.0
$ws.Range("L10").Value = 43.0
$ws.Range("M10").Value = 23.0
$ws.Range("N10").Value = 7.0
$ws.Range("O10").Value = 48.0
$ws.Range("P10").Value = 8.0
$ws.Range("Q10").Value = 31.0
$ws.Range("R10").Value = 1.0
$ws.Range("S10").Value = 29.0
$ws.Range("T10").Value = 21.0
$ws.Range("U10").Value = 34.0
$ws.Range("V10").Value = 30.0
$ws.Range("W10").Value = 16.0
$ws.Range("Y10").Value = 46.0
$ws.Range("Z10").Value = 38.0
$ws.Range("AA10").Value = 49.0
$ws.Range("AB10").Value = 39.0
$ws.Range("AC10").Value = 33.0
$ws.Range("AD10").Value = 45.0
$ws.Range("AE10").Value = 44.0
$ws.Range("AF10").Value = 19.0
$ws.Range("AG10").Value = 40.0
$ws.Range("AH10").Value = -1.0

# Row 11
$ws.Range("A11").Value = 0.0
$ws.Range("D11").Value = 469.0
$ws.Range("F11").Value = 4.0
$ws.Range("G11").Value = 14.0
$ws.Range("H11").Value = 6.0
$ws.Range("I11").Value = 27.0
$ws.Range("J11").Value = 51.0
$ws.Range("K11").Value = 32.0
$ws.Range("L11").Value = 26.0
$ws.Range("M11").Value = 22.0
$ws.Range("N11").Value = 28.0
$ws.Range("O11").Value = 3.0
$ws.Range("P11").Value = 36.0
$ws.Range("Q11").Value = 35.0
$ws.Range("R11").Value = 20.0
$ws.Range("S11").Value = 2.0
$ws.Range("T11").Value = 50.0
$ws.Range("U11").Value = 9.0
$ws.Range("V11").Value = 10.0
$ws.Range("W11").Value = 5.0
$ws.Range("X11").Value = 12.0
$ws.Range("Y11").Value = 17.0
$ws.Range("Z11").Value = 37.0
$ws.Range("AA11").Value = 15.0
$ws.Range("AB11").Value = 42.0
$ws.Range("AC11").Value = -1.0
$ws.Range("AD11").ClearContents()
$ws.Range("AE11").ClearContents()

# Row 12
$ws.Range("A12").Value = 0.0
$ws.Range("D12").Value = 468.0
$ws.Range("F12").Value = 41.0
$ws.Range("G12").Value = 13.0
$ws.Range("H12").Value = 47.0
$ws.Range("I12").Value = 18.0
$ws.Range("J12").Value = 25.0
$ws.Range("K12").Value = 24.0
$ws.Range("L12").Value = 43.0
$ws.Range("M12").Value = 23.0
$ws.Range("N12").Value = 7.0
$ws.Range("O12").Value = 48.0
$ws.Range("P12").Value = 8.0
$ws.Range("Q12").Value = 31.0
$ws.Range("R12").Value = 1.0
$ws.Range("S12").Value = 29.0
$ws.Range("T12").Value = 21.0
$ws.Range("U12").Value = 34.0
$ws.Range("V12").Value = 30.0
$ws.Range("W12").Value = 16.0
$ws.Range("Y12").Value = 46.0
$ws.Range("Z12").Value = 38.0
$ws.Range("AA12").Value = 49.0
$ws.Range("AB12").Value = 39.0
$ws.Range("AC12").Value = 33.0
$ws.Range("AD12").Value = 45.0
$ws.Range("AE12").Value = 44.0
$ws.Range("AF12").Value = 19.0
$ws.Range("AG12").Value = 40.0
$ws.Range("AH12").Value = -1.0
